$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2222657.8
$ws.Range("I6").Value = 4000148.8
$ws.Range("K6").Value = 12000446.4
$ws.Range("M6").Value = -12000334.4
$ws.Range("H29").Value = 4971
$ws.Range("J29").Value = 5879.9
$ws.Range("L29").Value = 17639.7
$ws.Range("N29").Value = -18201.7
$ws.Range("H74").Value = 9009.929
$ws.Range("I74").Value = 7266.875
$ws.Range("K74").Value = 7266.875
$ws.Range("M74").Value = -6330.875
$ws.Range("H76").Value = 12001.5
$ws.Range("I76").Value = 4999
$ws.Range("K76").Value = 4999
$ws.Range("M76").Value = -4684
$ws.Range("H77").Value = 9009.929
$ws.Range("I77").Value = 7266.875
$ws.Range("K77").Value = 36334.375
$ws.Range("M77").Value = -31654.375
$ws.Range("H79").Value = 12001.5
$ws.Range("I79").Value = 4999
$ws.Range("K79").Value = 4999
$ws.Range("M79").Value = -3907
$ws.Range("H116").Value = 17500.5
$ws.Range("I116").Value = 16499.75
$ws.Range("K116").Value = 16499.75
$ws.Range("M116").Value = -13057.75
$ws.Range("H132").Value = 1473.7059
$ws.Range("I132").Value = 1491.3939
$ws.Range("J132").Value = 890
$ws.Range("K132").Value = 4474.1817
$ws.Range("L132").Value = 2670
$ws.Range("M132").Value = -1944.1817
$ws.Range("N132").Value = -7730
$ws.Range("H138").Value = 3378.76
$ws.Range("I138").Value = 2876.6
$ws.Range("K138").Value = 8629.799999999999
$ws.Range("M138").Value = -3489.799999999999
$ws.Range("H141").Value = 1093.25
$ws.Range("I141").Value = 1093.25
$ws.Range("K141").Value = 3279.75
$ws.Range("M141").Value = 1900.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3599.932
$ws.Range("I61").Value = 3583.262
$ws.Range("K61").Value = 3583.262
$ws.Range("M61").Value = -3371.262
$ws.Range("H132").Value = 2939.913
$ws.Range("I132").Value = 2051.3171
$ws.Range("K132").Value = 6153.951300000001
$ws.Range("M132").Value = -3623.951300000001
$ws.Range("H136").Value = 3599.932
$ws.Range("I136").Value = 3583.262
$ws.Range("K136").Value = 10749.786
$ws.Range("M136").Value = -8199.786

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4996.25
$ws.Range("I86").Value = 4996.25
$ws.Range("K86").Value = 4996.25
$ws.Range("M86").Value = -3873.25
$ws.Range("H89").Value = 4996.25
$ws.Range("I89").Value = 4996.25
$ws.Range("K89").Value = 24981.25
$ws.Range("M89").Value = -19365.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 533.9286
$ws.Range("J7").Value = 381.5
$ws.Range("L7").Value = 381.5
$ws.Range("N7").Value = -607.5
$ws.Range("H8").Value = 2962.75
$ws.Range("I8").Value = 650
$ws.Range("J8").Value = 5275.5
$ws.Range("K8").Value = 650
$ws.Range("L8").Value = 5275.5
$ws.Range("M8").Value = -510
$ws.Range("N8").Value = -5555.5
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H86").Value = 17401.4
$ws.Range("J86").Value = 17401.4
$ws.Range("L86").Value = 17401.4
$ws.Range("N86").Value = -19647.4
$ws.Range("H89").Value = 17401.4
$ws.Range("J89").Value = 17401.4
$ws.Range("L89").Value = 87007
$ws.Range("N89").Value = -98239
$ws.Range("H107").Value = 1607.8948
$ws.Range("I107").Value = 1615.4445
$ws.Range("J107").Value = 1601.1
$ws.Range("K107").Value = 1615.4445
$ws.Range("L107").Value = 1601.1
$ws.Range("M107").Value = 304.5554999999999
$ws.Range("N107").Value = -5441.1
$ws.Range("H134").Value = 3761.8125
$ws.Range("I134").Value = 2438.8
$ws.Range("J134").Value = 5966.8335
$ws.Range("K134").Value = 7316.400000000001
$ws.Range("L134").Value = 17900.5005
$ws.Range("M134").Value = -4781.400000000001
$ws.Range("N134").Value = -22970.5005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 11666.667
$ws.Range("I3").Value = 3333.3333
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 9999.999899999999
$ws.Range("L3").Value = 60000
$ws.Range("M3").Value = -9887.999899999999
$ws.Range("N3").Value = -60224
$ws.Range("H49").Value = 125
$ws.Range("I49").Value = 125
$ws.Range("K49").Value = 375
$ws.Range("M49").Value = -219
$ws.Range("H56").Value = 6655
$ws.Range("I56").Value = 6655
$ws.Range("K56").Value = 6655
$ws.Range("M56").Value = -6125
$ws.Range("H80").Value = 8760.200000000001
$ws.Range("J80").Value = 10401
$ws.Range("L80").Value = 31203
$ws.Range("N80").Value = -33075
$ws.Range("H83").Value = 8760.200000000001
$ws.Range("J83").Value = 10401
$ws.Range("L83").Value = 93609
$ws.Range("N83").Value = -102969
$ws.Range("H132").Value = 4009.8333
$ws.Range("I132").Value = 2311
$ws.Range("K132").Value = 20799
$ws.Range("M132").Value = -18269

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H107").Value = 1963.091
$ws.Range("I107").Value = 825
$ws.Range("J107").Value = 2613.4285
$ws.Range("K107").Value = 825
$ws.Range("L107").Value = 2613.4285
$ws.Range("M107").Value = 1095
$ws.Range("N107").Value = -6453.4285
$ws.Range("H132").Value = 6368.25
$ws.Range("I132").Value = 3396.9375
$ws.Range("K132").Value = 10190.8125
$ws.Range("M132").Value = -7660.8125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8093.2383
$ws.Range("I7").Value = 4348.5
$ws.Range("J7").Value = 11497.546
$ws.Range("K7").Value = 4348.5
$ws.Range("L7").Value = 11497.546
$ws.Range("M7").Value = -4236.5
$ws.Range("N7").Value = -11721.546
$ws.Range("H46").Value = 4052
$ws.Range("J46").Value = 4777.3335
$ws.Range("L46").Value = 4777.3335
$ws.Range("N46").Value = -5153.3335
$ws.Range("H126").Value = 8093.2383
$ws.Range("I126").Value = 4348.5
$ws.Range("J126").Value = 11497.546
$ws.Range("K126").Value = 13045.5
$ws.Range("L126").Value = 34492.638
$ws.Range("M126").Value = -10575.5
$ws.Range("N126").Value = -39432.638
$ws.Range("H132").Value = 2959.9
$ws.Range("I132").Value = 1730.2693
$ws.Range("K132").Value = 5190.8079
$ws.Range("M132").Value = -2660.8079

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 9000
$ws.Range("I2").Value = 9000
$ws.Range("K2").Value = 9000
$ws.Range("M2").Value = -8888
$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 20000
$ws.Range("L20").Value = 20000
$ws.Range("N20").Value = -20480
$ws.Range("H49").Value = 29000
$ws.Range("J49").Value = 29000
$ws.Range("L49").Value = 29000
$ws.Range("N49").Value = -29460
$ws.Range("H81").Value = 4279.1333
$ws.Range("I81").Value = 2706.5386
$ws.Range("K81").Value = 5413.0772
$ws.Range("M81").Value = -4352.0772
$ws.Range("H84").Value = 4279.1333
$ws.Range("I84").Value = 2706.5386
$ws.Range("K84").Value = 27065.386
$ws.Range("M84").Value = -21761.386
$ws.Range("H132").Value = 5209.3193
$ws.Range("J132").Value = 14622.777
$ws.Range("L132").Value = 43868.331
$ws.Range("N132").Value = -48928.331
